$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing rows 138-145: columns P and AA get new values, AB = -P
$updates = @{
    138 = 2381
    139 = 2338
    140 = 2294
    141 = 2211
    142 = 2148
    143 = 2085
    144 = 2020
    145 = 1955
}

foreach ($row in $updates.Keys) {
    $val = $updates[$row]
    $ws.Range("P$row").Value = $val
    $ws.Range("AA$row").Value = $val
    $ws.Range("AB$row").Value = -$val
}

# Add new row 146 with the new period "01-09-2021"
$newRow = 146

# Write the period label as plain text without letting the smart-entry
# logic reinterpret it as a date serial (which would also mint a new
# number-format style). Compute it as a formula result, then paste back
# as a value so the stored cell is a plain shared string like its peers.
$ws.Range("A$newRow").Formula = '=TEXT("01-09-2021","@")'
$ws.Range("A$newRow").Copy()
$ws.Range("A$newRow").PasteSpecial(-4163)

$zeroCols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","Q","R","S","T","U","V","W","X","Y","Z")
foreach ($col in $zeroCols) {
    $ws.Range("$col$newRow").Value = 0
}

$ws.Range("P$newRow").Value = 1891
$ws.Range("AA$newRow").Value = 1891
$ws.Range("AB$newRow").Value = -1891
